$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 16: "Additional Effort" count changes from 3 to 5
$ws.Range("B16").Value = 5

# Copy formatting of the date cell A16 (style s="1") to the new date cells
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 17
$ws.Range("A17").Value = 41185
$ws.Range("B17").Value = 2.25
$ws.Range("D17").Value = "Implementation task overrun and stack usage, not tested yet"

# Row 18
$ws.Range("A18").Value = 41186
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = "Implementation application interrupts, not tested yet. Testing of setEvent/waitForEvent"

# Update the active selection to C18, matching the author's last edited cell
$ws.Range("C18").Select()
